$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '41.366.10'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '2.192.89'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.68'
$ws.Range("E5").Value = '  +5.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  +1.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.34'
$ws.Range("E7").Value = '  -1.09%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  +9.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.26'
$ws.Range("E10").Value = '  +5.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.70'
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0941'
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.13'
$ws.Range("E13").Value = '  +9.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '2.516.98'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.873'
$ws.Range("E16").Value = '  +5.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.56'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '2.185.00'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").Value = '41.236.29'
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").Value = '0.0₃0954'
$ws.Range("E20").Value = '  +2.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.26'
$ws.Range("E21").Value = '  +4.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.12'
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.10'
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("E24").Value = '  +3.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.87'
$ws.Range("E25").Value = '  +23.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.88'
$ws.Range("E26").Value = '  +7.88%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.51'
$ws.Range("E28").Value = '  +5.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.73'
$ws.Range("E30").Value = '  +1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.66'
$ws.Range("E31").Value = '  +3.18%  '
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.124'
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.49'
$ws.Range("E34").Value = '  +9.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0733'
$ws.Range("E35").Value = '  +5.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.52'
$ws.Range("E36").Value = '  +14.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.63'
$ws.Range("E37").Value = '  +2.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.11'
$ws.Range("E38").Value = '  +8.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0301'
$ws.Range("E39").Value = '  +13.61%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.78'
$ws.Range("E40").Value = '  +30.60%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.23'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.71'
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.77'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.203'
$ws.Range("E44").Value = '  +7.72%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.94'
$ws.Range("E45").Value = '  +5.84%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.69'
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("E47").Value = '  +4.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.14'
$ws.Range("E49").Value = '  +6.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("E50").Value = '  +2.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.32'
$ws.Range("E51").Value = '  +4.88%  '
